$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("A11").Value = 131167651
$ws.Range("B11").Value = 57884
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "NT"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = 100109
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "Tretåig hackspett"
$ws.Range("F11").NumberFormat = "General"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "Picoides tridactylus"
$ws.Range("G11").NumberFormat = "General"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("H11").NumberFormat = "General"
$ws.Range("I11").Value = $null
$ws.Range("J11").Value = $null
$ws.Range("K11").Value = $null
$ws.Range("L11").Value = $null
$ws.Range("M11").NumberFormat = "@"
$ws.Range("M11").Value = "äldre spår"
$ws.Range("M11").NumberFormat = "General"
$ws.Range("N11").Value = $null
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value = "Vallsjöberget, Ång"
$ws.Range("P11").NumberFormat = "General"
$ws.Range("Q11").Value = 613285
$ws.Range("R11").Value = 6997537
$ws.Range("S11").Value = 10
$ws.Range("T11").NumberFormat = "@"
$ws.Range("T11").Value = "Västernorrland"
$ws.Range("T11").NumberFormat = "General"
$ws.Range("U11").NumberFormat = "@"
$ws.Range("U11").Value = "Sollefteå"
$ws.Range("U11").NumberFormat = "General"
$ws.Range("V11").NumberFormat = "@"
$ws.Range("V11").Value = "Ångermanland"
$ws.Range("V11").NumberFormat = "General"
$ws.Range("W11").NumberFormat = "@"
$ws.Range("W11").Value = "Sollefteå"
$ws.Range("W11").NumberFormat = "General"
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("Y11").Value = "2026-02-14"
$ws.Range("Y11").NumberFormat = "General"
$ws.Range("AA11").NumberFormat = "@"
$ws.Range("AA11").Value = "2026-02-14"
$ws.Range("AA11").NumberFormat = "General"
$ws.Range("AC11").NumberFormat = "@"
$ws.Range("AC11").Value = "Äldre ringhack på tall"
$ws.Range("AC11").NumberFormat = "General"
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AF11").Value = $null
$ws.Range("AG11").Value = $false
$ws.Range("AJ11").Value = $null
$ws.Range("AK11").Value = $null
$ws.Range("AO11").Value = $null
$ws.Range("AT11").Value = $null
$ws.Range("AW11").NumberFormat = "@"
$ws.Range("AW11").Value = "Markus Borja"
$ws.Range("AW11").NumberFormat = "General"
$ws.Range("AX11").NumberFormat = "@"
$ws.Range("AX11").Value = "Markus Borja"
$ws.Range("AX11").NumberFormat = "General"
$ws.Range("AY11").Value = $null

# Row 13
$ws.Range("A13").Value = 131167669
$ws.Range("B13").Value = 79244
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "NT"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = 6425
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "Garnlav"
$ws.Range("F13").NumberFormat = "General"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "Alectoria sarmentosa"
$ws.Range("G13").NumberFormat = "General"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "(Ach.) Ach."
$ws.Range("H13").NumberFormat = "General"
$ws.Range("I13").Value = $null
$ws.Range("J13").Value = $null
$ws.Range("K13").Value = $null
$ws.Range("L13").Value = $null
$ws.Range("M13").Value = $null
$ws.Range("N13").Value = $null
$ws.Range("P13").NumberFormat = "@"
$ws.Range("P13").Value = "Vallsjöberget, Ång"
$ws.Range("P13").NumberFormat = "General"
$ws.Range("Q13").Value = 613256
$ws.Range("R13").Value = 6997380
$ws.Range("S13").Value = 10
$ws.Range("T13").NumberFormat = "@"
$ws.Range("T13").Value = "Västernorrland"
$ws.Range("T13").NumberFormat = "General"
$ws.Range("U13").NumberFormat = "@"
$ws.Range("U13").Value = "Sollefteå"
$ws.Range("U13").NumberFormat = "General"
$ws.Range("V13").NumberFormat = "@"
$ws.Range("V13").Value = "Ångermanland"
$ws.Range("V13").NumberFormat = "General"
$ws.Range("W13").NumberFormat = "@"
$ws.Range("W13").Value = "Sollefteå"
$ws.Range("W13").NumberFormat = "General"
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = "2026-02-14"
$ws.Range("Y13").NumberFormat = "General"
$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value = "2026-02-14"
$ws.Range("AA13").NumberFormat = "General"
$ws.Range("AC13").Value = $null
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AF13").Value = $null
$ws.Range("AG13").Value = $false
$ws.Range("AJ13").NumberFormat = "@"
$ws.Range("AJ13").Value = "tall"
$ws.Range("AJ13").NumberFormat = "General"
$ws.Range("AK13").NumberFormat = "@"
$ws.Range("AK13").Value = "Pinus sylvestris"
$ws.Range("AK13").NumberFormat = "General"
$ws.Range("AO13").NumberFormat = "@"
$ws.Range("AO13").Value = "Pinus sylvestris"
$ws.Range("AO13").NumberFormat = "General"
$ws.Range("AT13").Value = $null
$ws.Range("AW13").NumberFormat = "@"
$ws.Range("AW13").Value = "Markus Borja"
$ws.Range("AW13").NumberFormat = "General"
$ws.Range("AX13").NumberFormat = "@"
$ws.Range("AX13").Value = "Markus Borja"
$ws.Range("AX13").NumberFormat = "General"
$ws.Range("AY13").Value = $null

# Row 16
$ws.Range("A16").Value = 131167655
$ws.Range("B16").Value = 57884
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "NT"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = 100109
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "Tretåig hackspett"
$ws.Range("F16").NumberFormat = "General"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "Picoides tridactylus"
$ws.Range("G16").NumberFormat = "General"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "(Linnaeus, 1758)"
$ws.Range("H16").NumberFormat = "General"
$ws.Range("I16").Value = $null
$ws.Range("K16").Value = $null
$ws.Range("L16").Value = $null
$ws.Range("M16").NumberFormat = "@"
$ws.Range("M16").Value = "färska spår"
$ws.Range("M16").NumberFormat = "General"
$ws.Range("N16").Value = $null
$ws.Range("P16").NumberFormat = "@"
$ws.Range("P16").Value = "Vallsjöberget, Ång"
$ws.Range("P16").NumberFormat = "General"
$ws.Range("Q16").Value = 613285
$ws.Range("R16").Value = 6997398
$ws.Range("S16").Value = 10
$ws.Range("T16").NumberFormat = "@"
$ws.Range("T16").Value = "Västernorrland"
$ws.Range("T16").NumberFormat = "General"
$ws.Range("U16").NumberFormat = "@"
$ws.Range("U16").Value = "Sollefteå"
$ws.Range("U16").NumberFormat = "General"
$ws.Range("V16").NumberFormat = "@"
$ws.Range("V16").Value = "Ångermanland"
$ws.Range("V16").NumberFormat = "General"
$ws.Range("W16").NumberFormat = "@"
$ws.Range("W16").Value = "Sollefteå"
$ws.Range("W16").NumberFormat = "General"
$ws.Range("Y16").NumberFormat = "@"
$ws.Range("Y16").Value = "2026-02-14"
$ws.Range("Y16").NumberFormat = "General"
$ws.Range("AA16").NumberFormat = "@"
$ws.Range("AA16").Value = "2026-02-14"
$ws.Range("AA16").NumberFormat = "General"
$ws.Range("AC16").NumberFormat = "@"
$ws.Range("AC16").Value = "Färska ringhack på tall"
$ws.Range("AC16").NumberFormat = "General"
$ws.Range("AD16").Value = $false
$ws.Range("AE16").Value = $false
$ws.Range("AG16").Value = $false
$ws.Range("AT16").Value = $null
$ws.Range("AW16").NumberFormat = "@"
$ws.Range("AW16").Value = "Markus Borja"
$ws.Range("AW16").NumberFormat = "General"
$ws.Range("AX16").NumberFormat = "@"
$ws.Range("AX16").Value = "Markus Borja"
$ws.Range("AX16").NumberFormat = "General"
$ws.Range("AY16").Value = $null

# Row 17
$ws.Range("A17").Value = 131167658
$ws.Range("B17").Value = 57073
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "LC"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = 100138
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "Tjäder"
$ws.Range("F17").NumberFormat = "General"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "Tetrao urogallus"
$ws.Range("G17").NumberFormat = "General"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "Linnaeus, 1758"
$ws.Range("H17").NumberFormat = "General"
$ws.Range("I17").Value = $null
$ws.Range("K17").Value = $null
$ws.Range("L17").Value = $null
$ws.Range("M17").NumberFormat = "@"
$ws.Range("M17").Value = "färsk spillning"
$ws.Range("M17").NumberFormat = "General"
$ws.Range("N17").Value = $null
$ws.Range("P17").NumberFormat = "@"
$ws.Range("P17").Value = "Vallsjöberget, Ång"
$ws.Range("P17").NumberFormat = "General"
$ws.Range("Q17").Value = 613330
$ws.Range("R17").Value = 6997326
$ws.Range("S17").Value = 10
$ws.Range("T17").NumberFormat = "@"
$ws.Range("T17").Value = "Västernorrland"
$ws.Range("T17").NumberFormat = "General"
$ws.Range("U17").NumberFormat = "@"
$ws.Range("U17").Value = "Sollefteå"
$ws.Range("U17").NumberFormat = "General"
$ws.Range("V17").NumberFormat = "@"
$ws.Range("V17").Value = "Ångermanland"
$ws.Range("V17").NumberFormat = "General"
$ws.Range("W17").NumberFormat = "@"
$ws.Range("W17").Value = "Sollefteå"
$ws.Range("W17").NumberFormat = "General"
$ws.Range("Y17").NumberFormat = "@"
$ws.Range("Y17").Value = "2026-02-14"
$ws.Range("Y17").NumberFormat = "General"
$ws.Range("AA17").NumberFormat = "@"
$ws.Range("AA17").Value = "2026-02-14"
$ws.Range("AA17").NumberFormat = "General"
$ws.Range("AC17").Value = $null
$ws.Range("AD17").Value = $false
$ws.Range("AE17").Value = $false
$ws.Range("AG17").Value = $false
$ws.Range("AT17").Value = $null
$ws.Range("AW17").NumberFormat = "@"
$ws.Range("AW17").Value = "Markus Borja"
$ws.Range("AW17").NumberFormat = "General"
$ws.Range("AX17").NumberFormat = "@"
$ws.Range("AX17").Value = "Markus Borja"
$ws.Range("AX17").NumberFormat = "General"
$ws.Range("AY17").Value = $null

# Row 19
$ws.Range("A19").Value = 131167652
$ws.Range("B19").Value = 57884
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "NT"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = 100109
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "Tretåig hackspett"
$ws.Range("F19").NumberFormat = "General"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "Picoides tridactylus"
$ws.Range("G19").NumberFormat = "General"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "(Linnaeus, 1758)"
$ws.Range("H19").NumberFormat = "General"
$ws.Range("I19").NumberFormat = "@"
$ws.Range("I19").Value = "1"
$ws.Range("I19").NumberFormat = "General"
$ws.Range("K19").Value = $null
$ws.Range("L19").Value = $null
$ws.Range("M19").NumberFormat = "@"
$ws.Range("M19").Value = "födosökande"
$ws.Range("M19").NumberFormat = "General"
$ws.Range("N19").NumberFormat = "@"
$ws.Range("N19").Value = "observerad"
$ws.Range("N19").NumberFormat = "General"
$ws.Range("P19").NumberFormat = "@"
$ws.Range("P19").Value = "Vallsjöberget, Ång"
$ws.Range("P19").NumberFormat = "General"
$ws.Range("Q19").Value = 613264
$ws.Range("R19").Value = 6997532
$ws.Range("S19").Value = 10
$ws.Range("T19").NumberFormat = "@"
$ws.Range("T19").Value = "Västernorrland"
$ws.Range("T19").NumberFormat = "General"
$ws.Range("U19").NumberFormat = "@"
$ws.Range("U19").Value = "Sollefteå"
$ws.Range("U19").NumberFormat = "General"
$ws.Range("V19").NumberFormat = "@"
$ws.Range("V19").Value = "Ångermanland"
$ws.Range("V19").NumberFormat = "General"
$ws.Range("W19").NumberFormat = "@"
$ws.Range("W19").Value = "Sollefteå"
$ws.Range("W19").NumberFormat = "General"
$ws.Range("Y19").NumberFormat = "@"
$ws.Range("Y19").Value = "2026-02-14"
$ws.Range("Y19").NumberFormat = "General"
$ws.Range("Z19").NumberFormat = "@"
$ws.Range("Z19").Value = "12:46"
$ws.Range("Z19").NumberFormat = "General"
$ws.Range("AA19").NumberFormat = "@"
$ws.Range("AA19").Value = "2026-02-14"
$ws.Range("AA19").NumberFormat = "General"
$ws.Range("AB19").NumberFormat = "@"
$ws.Range("AB19").Value = "12:48"
$ws.Range("AB19").NumberFormat = "General"
$ws.Range("AC19").NumberFormat = "@"
$ws.Range("AC19").Value = "Hackspetten syns i profil på ett smalt träd i centrum av den tagna bilden."
$ws.Range("AC19").NumberFormat = "General"
$ws.Range("AD19").Value = $false
$ws.Range("AE19").Value = $false
$ws.Range("AG19").Value = $false
$ws.Range("AT19").Value = $null
$ws.Range("AW19").NumberFormat = "@"
$ws.Range("AW19").Value = "Markus Borja"
$ws.Range("AW19").NumberFormat = "General"
$ws.Range("AX19").NumberFormat = "@"
$ws.Range("AX19").Value = "Markus Borja"
$ws.Range("AX19").NumberFormat = "General"
$ws.Range("AY19").Value = $null

# Row 20
$ws.Range("A20").Value = 131167654
$ws.Range("B20").Value = 57884
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "NT"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = 100109
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "Tretåig hackspett"
$ws.Range("F20").NumberFormat = "General"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "Picoides tridactylus"
$ws.Range("G20").NumberFormat = "General"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "(Linnaeus, 1758)"
$ws.Range("H20").NumberFormat = "General"
$ws.Range("I20").Value = $null
$ws.Range("K20").Value = $null
$ws.Range("L20").Value = $null
$ws.Range("M20").NumberFormat = "@"
$ws.Range("M20").Value = "äldre spår"
$ws.Range("M20").NumberFormat = "General"
$ws.Range("N20").Value = $null
$ws.Range("P20").NumberFormat = "@"
$ws.Range("P20").Value = "Vallsjöberget, Ång"
$ws.Range("P20").NumberFormat = "General"
$ws.Range("Q20").Value = 613254
$ws.Range("R20").Value = 6997565
$ws.Range("S20").Value = 10
$ws.Range("T20").NumberFormat = "@"
$ws.Range("T20").Value = "Västernorrland"
$ws.Range("T20").NumberFormat = "General"
$ws.Range("U20").NumberFormat = "@"
$ws.Range("U20").Value = "Sollefteå"
$ws.Range("U20").NumberFormat = "General"
$ws.Range("V20").NumberFormat = "@"
$ws.Range("V20").Value = "Ångermanland"
$ws.Range("V20").NumberFormat = "General"
$ws.Range("W20").NumberFormat = "@"
$ws.Range("W20").Value = "Sollefteå"
$ws.Range("W20").NumberFormat = "General"
$ws.Range("Y20").NumberFormat = "@"
$ws.Range("Y20").Value = "2026-02-14"
$ws.Range("Y20").NumberFormat = "General"
$ws.Range("Z20").Value = $null
$ws.Range("AA20").NumberFormat = "@"
$ws.Range("AA20").Value = "2026-02-14"
$ws.Range("AA20").NumberFormat = "General"
$ws.Range("AB20").Value = $null
$ws.Range("AC20").NumberFormat = "@"
$ws.Range("AC20").Value = "Äldre ringhack på tall"
$ws.Range("AC20").NumberFormat = "General"
$ws.Range("AD20").Value = $false
$ws.Range("AE20").Value = $false
$ws.Range("AG20").Value = $false
$ws.Range("AT20").Value = $null
$ws.Range("AW20").NumberFormat = "@"
$ws.Range("AW20").Value = "Markus Borja"
$ws.Range("AW20").NumberFormat = "General"
$ws.Range("AX20").NumberFormat = "@"
$ws.Range("AX20").Value = "Markus Borja"
$ws.Range("AX20").NumberFormat = "General"
$ws.Range("AY20").Value = $null
